# Fix template spelling, dicipline -> discipline in vars_meta_data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vars_meta_data")

$ws.Range("I1").Value = "var_discipline"
$ws.Range("I2").Value = "< associated discipline(s) (<100 chars) (examples: [Physics, Chemistry, Biology, BioGeoChemistry, etc..] )>  ↓"

# The saved file also shows this sheet became the active tab/selection
# (reflecting the state of the workbook when the author saved the fix).
$ws.Activate() | Out-Null
$ws.Range("I2").Select() | Out-Null
